$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.743.63"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.077.18"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'233.78"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'58.18"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'0.106"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'14.89"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.383.33"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "'20.91"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "'0.774"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "2.053.32"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "37.703.27"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'71.14"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D22").Value = "'227.83"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "'169.41"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").Value = "'0.139"
$ws.Range("E27").Value = "  +3.92%  "
$ws.Range("D28").Value = "'9.00"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "'19.45"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("D32").Value = "'4.68"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "'0.0630"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").Value = "'4.66"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("D36").Value = "'1.83"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("D37").Value = "'3.40"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "'5.36"
$ws.Range("E39").Value = "  -4.58%  "
$ws.Range("D40").Value = "'0.0980"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").Value = "'98.13"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "1.453.67"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").Value = "'16.61"
$ws.Range("E45").Value = "  +6.58%  "
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "'4.25"
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").Value = "'7.38"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").Value = "'3.02"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "2.267.13"
$ws.Range("E51").Value = "  -1.38%  "
